$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 3.2
$ws.Range("I4").Value = 2.25
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 9.5

$ws.Range("G5").Value = 2.05
$ws.Range("S5").Value = 1.67

$ws.Range("G6").Value = 3.5
$ws.Range("I6").Value = 2.05
$ws.Range("J6").Value = 1.07
$ws.Range("K6").Value = 9
$ws.Range("R6").Value = 1.83
$ws.Range("S6").Value = 1.83
$ws.Range("T6").Value = 10
$ws.Range("V6").Value = 13
$ws.Range("W6").Value = 41
$ws.Range("Y6").Value = 41
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 9.5
$ws.Range("AH6").Value = 17
$ws.Range("AJ6").Value = 301

$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 2.77
$ws.Range("K7").Value = 4.35
$ws.Range("M7").Value = 2.02
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 2.32
$ws.Range("T7").Value = 6
$ws.Range("U7").Value = 13.5
$ws.Range("V7").Value = 12
$ws.Range("W7").Value = 40
$ws.Range("X7").Value = 37
$ws.Range("Z7").Value = 4.35
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 175
$ws.Range("AD7").Value = 5.6
$ws.Range("AE7").Value = 11.75
$ws.Range("AF7").Value = 11.5
$ws.Range("AG7").Value = 35
$ws.Range("AH7").Value = 35
$ws.Range("AI7").Value = 65

$ws.Range("G9").Value = 1.7
$ws.Range("L9").Value = 1.3
$ws.Range("M9").Value = 3.4
$ws.Range("R9").Value = 1.91
$ws.Range("S9").Value = 1.8

$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 6.5
$ws.Range("L10").Value = 1.23
$ws.Range("M10").Value = 3.4
$ws.Range("N10").Value = 1.7
$ws.Range("O10").Value = 1.93
$ws.Range("R10").Value = 1.83
$ws.Range("S10").Value = 1.78
$ws.Range("U10").Value = 6.9
$ws.Range("V10").Value = 8
$ws.Range("W10").Value = 9.75
$ws.Range("X10").Value = 11.5
$ws.Range("Z10").Value = 11.5
$ws.Range("AA10").Value = 8.25
$ws.Range("AB10").Value = 18
$ws.Range("AD10").Value = 17.5
$ws.Range("AE10").Value = 45
$ws.Range("AF10").Value = 20
$ws.Range("AG10").Value = 150
$ws.Range("AH10").Value = 70
$ws.Range("AI10").Value = 65
$ws.Range("AJ10").Value = 700

$ws.Range("G11").Value = 2.62
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 2.47
$ws.Range("M11").Value = 4.4
$ws.Range("R11").Value = 1.47
$ws.Range("S11").Value = 2.35
$ws.Range("T11").Value = 12.5
$ws.Range("V11").Value = 9.75
$ws.Range("W11").Value = 32
$ws.Range("X11").Value = 19
$ws.Range("Y11").Value = 21
$ws.Range("AA11").Value = 6.8
$ws.Range("AD11").Value = 10.75
$ws.Range("AF11").Value = 9.25
$ws.Range("AG11").Value = 28
$ws.Range("AH11").Value = 18.5
$ws.Range("AI11").Value = 22

$ws.Range("L12").Value = 1.29
$ws.Range("M12").Value = 3.3
$ws.Range("R12").Value = 1.73
$ws.Range("S12").Value = 2

$ws.Range("S13").Value = 2.11

$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 2.63

$ws.Range("R16").Value = 1.67
$ws.Range("S16").Value = 2.1
$ws.Range("T16").Value = 9

$ws.Range("G17").Value = 1.48
$ws.Range("H17").Value = 4.5
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 1.03
$ws.Range("K17").Value = 17
$ws.Range("N17").Value = 1.53
$ws.Range("O17").Value = 2.4
$ws.Range("P17").Value = 1.29
$ws.Range("Q17").Value = 3.5
$ws.Range("R17").Value = 1.67
$ws.Range("AE17").Value = 34

$ws.Range("S18").Value = 1.67

$ws.Range("N20").Value = 1.9
$ws.Range("O20").Value = 1.9

$ws.Range("O22").Value = 1.85

$ws.Range("G23").Value = 1.45
$ws.Range("I23").Value = 5.25
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 1.03
$ws.Range("N23").Value = 1.5
$ws.Range("O23").Value = 2.5
$ws.Range("R23").Value = 1.62
$ws.Range("S23").Value = 2.2
$ws.Range("T23").Value = 10
$ws.Range("U23").Value = 9
$ws.Range("V23").Value = 8.5
$ws.Range("W23").Value = 11
$ws.Range("Y23").Value = 21
$ws.Range("Z23").Value = 19
$ws.Range("AD23").Value = 19
$ws.Range("AH23").Value = 41
$ws.Range("AJ23").Value = 151

$ws.Range("G24").Value = 1.8
$ws.Range("I24").Value = 3.5
$ws.Range("J24").Value = 1.05
$ws.Range("K24").Value = 8.5
$ws.Range("L24").Value = 1.29
$ws.Range("M24").Value = 3.5
$ws.Range("N24").Value = 1.88
$ws.Range("O24").Value = 1.93
$ws.Range("U24").Value = 9

$ws.Range("AJ25").Value = 500

$ws.Range("J26").Value = 1.08
$ws.Range("K26").Value = 6.4
$ws.Range("M26").Value = 2.75
$ws.Range("N26").Value = 2.15
$ws.Range("O26").Value = 1.62
$ws.Range("Q26").Value = 2.57
$ws.Range("R26").Value = 1.91
$ws.Range("S26").Value = 1.8
$ws.Range("V26").Value = 11.75
$ws.Range("Z26").Value = 6.4
$ws.Range("AB26").Value = 16
$ws.Range("AD26").Value = 6.6
$ws.Range("AE26").Value = 9.75
$ws.Range("AF26").Value = 9
$ws.Range("AH26").Value = 19
$ws.Range("AJ26").Value = 800

$ws.Range("G28").Value = 1.6
$ws.Range("H28").Value = 3.8
$ws.Range("I28").Value = 5.3
$ws.Range("J28").Value = 1.03
$ws.Range("K28").Value = 9
$ws.Range("L28").Value = 1.17
$ws.Range("M28").Value = 4.4
$ws.Range("N28").Value = 1.53
$ws.Range("O28").Value = 2.32
$ws.Range("P28").Value = 1.3
$ws.Range("T28").Value = 9.25
$ws.Range("U28").Value = 9.25
$ws.Range("V28").Value = 7.8
$ws.Range("W28").Value = 13
$ws.Range("X28").Value = 11.25
$ws.Range("Z28").Value = 9
$ws.Range("AA28").Value = 7.8
$ws.Range("AD28").Value = 19
$ws.Range("AE28").Value = 37
$ws.Range("AF28").Value = 16.5
$ws.Range("AG28").Value = 110

$ws.Range("G29").Value = 2.05
$ws.Range("I29").Value = 3.6
$ws.Range("R29").Value = 1.85
$ws.Range("S29").Value = 1.75
$ws.Range("T29").Value = 6.4
$ws.Range("U29").Value = 9.25
$ws.Range("W29").Value = 19
$ws.Range("X29").Value = 18
$ws.Range("AA29").Value = 6.1
$ws.Range("AD29").Value = 9.25
$ws.Range("AE29").Value = 18.5
$ws.Range("AF29").Value = 12.5
$ws.Range("AG29").Value = 55
$ws.Range("AH29").Value = 37
$ws.Range("AJ29").Value = 800
